$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A1 text (was: &; &amp; &quot; &lt; &gt; &apos;  -> now: &&apos; &amp; &quot; &lt; &gt; &apos;)
$ws.Range("A1").Value = '&&apos; &amp; &quot; &lt; &gt; &apos;'

# Rename the sheet (was the long escaped name, now: & & " > < )
$ws.Name = '& & " > < '
